# Auto-generated edit script for cryptos.xlsx update
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = '28.308.65'
$ws.Range("E2").Value = '  -2.27%  '
$ws.Range("D3").Value = '1.872.11'
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = '318.63'
$ws.Range("D5").Style = "Normal"
$ws.Range("E5").Value = '  -1.84%  '
$ws.Range("E6").Value = '  -0.08%  '
$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = '0.4402'
$ws.Range("D7").Style = "Normal"
$ws.Range("E7").Value = '  -4.16%  '
$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = '0.3702'
$ws.Range("D8").Style = "Normal"
$ws.Range("E8").Value = '  -3.22%  '
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = '0.07517'
$ws.Range("D9").Style = "Normal"
$ws.Range("E9").Value = '  -2.52%  '
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = '0.9396'
$ws.Range("D10").Style = "Normal"
$ws.Range("E10").Value = '  -4.06%  '
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = '21.40'
$ws.Range("D11").Style = "Normal"
$ws.Range("E11").Value = '  -2.86%  '
$ws.Range("D12").Value = '1.892.68'
$ws.Range("E12").Value = '  -1.57%  '
$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = '6.731'
$ws.Range("D13").Style = "Normal"
$ws.Range("E13").Value = '  -2.94%  '
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = '5.460'
$ws.Range("D14").Style = "Normal"
$ws.Range("E14").Value = '  -3.53%  '
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = '0.06862'
$ws.Range("D15").Style = "Normal"
$ws.Range("E15").Value = '  -2.46%  '
$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = '1.003'
$ws.Range("D16").Style = "Normal"
$ws.Range("E16").Value = '  -0.07%  '
$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = '82.30'
$ws.Range("D17").Style = "Normal"
$ws.Range("E17").Value = '  -1.68%  '
$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = '0.000009045'
$ws.Range("D18").Style = "Normal"
$ws.Range("E18").Value = '  -4.66%  '
$ws.Range("E20").Value = '  -4.16%  '
$ws.Range("D21").Value = '28.304.52'
$ws.Range("E21").Value = '  -2.27%  '
$ws.Range("E22").Value = '  -3.16%  '
$ws.Range("E23").Value = '  -0.02%  '
$ws.Range("D24").Value = '2.127.78'
$ws.Range("E24").Value = '  -1.31%  '
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = '2.031'
$ws.Range("D25").Style = "Normal"
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = '155.07'
$ws.Range("D26").Style = "Normal"
$ws.Range("E26").Value = '  -1.87%  '
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = '18.37'
$ws.Range("D27").Style = "Normal"
$ws.Range("E27").Value = '  -3.72%  '
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = '5.341'
$ws.Range("D28").Style = "Normal"
$ws.Range("E28").Value = '  -5.34%  '
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = '113.87'
$ws.Range("D29").Style = "Normal"
$ws.Range("E29").Value = '  -2.95%  '
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = '1.730'
$ws.Range("D30").Style = "Normal"
$ws.Range("E30").Value = '  -6.54%  '
$ws.Range("E31").Value = '  -2.30%  '
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = '0.7996'
$ws.Range("D32").Style = "Normal"
$ws.Range("E32").Value = '  -7.65%  '
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = '4.850'
$ws.Range("D33").Style = "Normal"
$ws.Range("E33").Value = '  -4.47%  '
$ws.Range("E34").Value = '  -5.56%  '
$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = '2.918'
$ws.Range("D35").Style = "Normal"
$ws.Range("E35").Value = '  -1.23%  '
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = '1.002'
$ws.Range("D36").Style = "Normal"
$ws.Range("E36").Value = '  -0.03%  '
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = '1.129'
$ws.Range("D37").Style = "Normal"
$ws.Range("E37").Value = '  -1.56%  '
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = '0.05449'
$ws.Range("D38").Style = "Normal"
$ws.Range("E38").Value = '  -4.71%  '
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = '0.01960'
$ws.Range("D39").Style = "Normal"
$ws.Range("E39").Value = '  -3.77%  '
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = '2.942'
$ws.Range("D40").Style = "Normal"
$ws.Range("E40").Value = '  +5.85%  '
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = '7.115'
$ws.Range("D41").Style = "Normal"
$ws.Range("E41").Value = '  -3.89%  '
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = '0.5263'
$ws.Range("D42").Style = "Normal"
$ws.Range("E42").Value = '  -4.31%  '
$ws.Range("E43").Value = '  -3.98%  '
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = '8.775'
$ws.Range("D44").Style = "Normal"
$ws.Range("E44").Value = '  -5.78%  '
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = '0.06759'
$ws.Range("D45").Style = "Normal"
$ws.Range("E45").Value = '  -1.11%  '
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = '0.4878'
$ws.Range("D46").Style = "Normal"
$ws.Range("E46").Value = '  -5.73%  '
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = '1.984'
$ws.Range("D47").Style = "Normal"
$ws.Range("E47").Value = '  -3.77%  '
$ws.Range("B48").Value = 'Quant'
$ws.Range("C48").Value = 'https://coinranking.com/coin/bauj_21eYVwso+quant-qnt'
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = '107.87'
$ws.Range("D48").Style = "Normal"
$ws.Range("E48").Value = '  -2.48%  '
$ws.Range("B49").Value = 'EnergySwap'
$ws.Range("C49").Value = 'https://coinranking.com/coin/SbWqqTui-+energyswap-ens'
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = '10.52'
$ws.Range("D49").Style = "Normal"
$ws.Range("E49").Value = '  -6.73%  '
$ws.Range("E50").Value = '  -5.41%  '
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = '1.682'
$ws.Range("D51").Style = "Normal"
$ws.Range("E51").Value = '  -5.17%  '
